# GNB2020_bgt_revs.xlsx bugfix:
# Four municipalities (Bathurst, Campbellton, Dieppe, Edmundston) were
# incorrectly excluded from the 2020 revenue data. This inserts four new
# rows right after the header (in alphabetical order, matching the rest
# of the table) and populates them with the correct revenue figures,
# then grows the worksheet table (ListObject) to cover the new range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert four blank rows right below the header row, pushing the
# existing municipalities (Fredericton...Tracy) down by four rows.
$ws.Rows("2:5").Insert()

# The inserted rows picked up the header's formatting; restore the
# normal data-row formatting by copying it down from row 6 (which now
# holds what used to be row 2 before the insert).
$srcRow = $ws.Range("A6:J6")
$destRows = $ws.Range("A2:J5")
$srcRow.Copy()
$destRows.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the four newly-inserted rows with the correct data.
$newRows = @(
  @{ Row = 2; Name = "Bathurst";    Values = @(19153436,3624702,1013496,1164400,472498,0,500000,0,25928532) },
  @{ Row = 3; Name = "Campbellton"; Values = @(10686193,1889297,156231,1093136,503816,46500,405200,66867,14847240) },
  @{ Row = 4; Name = "Dieppe";      Values = @(53323334,1095459,345000,1697000,978500,0,926441,583007,58948741) },
  @{ Row = 5; Name = "Edmundston";  Values = @(23720297,5564820,1718846,1342789,725500,5000,2420632,12801,35510685) }
)

foreach ($entry in $newRows) {
  $r = $entry.Row
  $ws.Cells.Item($r, 1).Value = $entry.Name
  $vals = $entry.Values
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
  }
}

# Grow the table (ListObject) so it covers the four new rows too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J96"))
